$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook calculation properties -------------------------------------
# Source diff turns on single-threaded (non concurrent) calculation
# (<calcPr .. concurrentCalc="0"/>). Mirror that through the Application's
# MultiThreadedCalculation object.
try {
    $excel.MultiThreadedCalculation.Enabled = $false
} catch {
}

# --- New header columns: "Level Id" (H) and "Parent Id" (I) --------------
# Row 4 holds the column headers for the "Address" table; add two more
# headers with the same bold style as the existing ones (G4, "AutoNumber").
$ws.Range("H4").Value = "Level Id"
$ws.Range("I4").Value = "Parent Id"
$ws.Range("H4:I4").Font.Bold = $true

# --- Move the LVL-index values out of column D into the new column H -----
# (and clear the now-unused D/E cells, keeping their existing yellow-fill
# style intact).

# Row 8 ("Address" / LVL row): D8 had the level index (1).
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("H8").Value = 1

# Row 11 ("Properties" / LVL row): D11 had the level index (2),
# E11 had the parent level index (1).
$ws.Range("D11").ClearContents()
$ws.Range("E11").ClearContents()
$ws.Range("H11").Value = 2
$ws.Range("I11").Value = 1

# Row 14 ("Features" / LVL row): D14 had the level index (3).
$ws.Range("D14").ClearContents()
$ws.Range("E14").ClearContents()
$ws.Range("H14").Value = 3

# --- Selection: the author ended up with I13 selected --------------------
$ws.Range("I13").Select() | Out-Null
